$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Simple single-run paragraph replacements (safe with Find & Replace).
Replace-Text "Currículum Vitae - Patti Fernandez" "Currículum Vitae: Patti Fernandez"
Replace-Text "ABC Studios: Animator principal (enero de 2018 - Presente)" "ABC Studios: Animadora principal (enero de 2018 - actualidad)"
Replace-Text "XYZ Media: Animator Senior (jun 2015 - dic 2017)" "XYZ Media: Animadora principal (junio de 2015 - diciembre de 2017)"
Replace-Text "MNO Entertainment: Junior Animator (Sep 2012 - Mayo de 2015)" "MNO Entertainment: Animadora júnior (septiembre de 2012 - mayo de 2015)"

# The "El arte de la animación..." run sits between two other runs that
# share the exact same run formatting. A plain Find/Replace (or a direct
# Range.Text assignment) on just that run causes this runtime to coalesce
# it with its neighboring same-formatted runs, which would also swallow
# the unrelated " " and "Nueva York: Penguin Books." runs that must stay
# untouched. Briefly nudging a run's formatting away from its neighbors
# before/after the text edit keeps run boundaries intact; we restore the
# original size afterward so no visible/semantic formatting changes.
$target = $d.Content.Duplicate
$foundTarget = $target.Find.Execute("El arte de la animación: una guía para principiantes.")
if ($foundTarget) {
    $originalSize = $target.Font.Size

    # Guard the trailing boundary too: touching the target run's text can
    # otherwise cascade a merge of the following " " + "Nueva York: Penguin
    # Books." runs even though neither is being edited.
    $afterRun = $d.Content.Duplicate
    $foundAfter = $afterRun.Find.Execute("Nueva York: Penguin Books.")
    if ($foundAfter) {
        $boundary = $d.Range($afterRun.Start - 1, $afterRun.Start)
        $boundarySize = $boundary.Font.Size
        $boundary.Font.Size = $boundarySize + 1
    }

    $target.Font.Size = $originalSize + 1
    $target.Text = "The Art of Animation: A Guide for Beginners."

    $newTarget = $d.Range($target.Start, $target.Start + "The Art of Animation: A Guide for Beginners.".Length)
    $newTarget.Font.Size = $originalSize

    if ($foundAfter) {
        $boundary2 = $d.Range($newTarget.End, $newTarget.End + 1)
        $boundary2.Font.Size = $boundarySize
    }
}
